$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 10207
$ws1.Range("F5").Value = 737
$ws1.Range("F6").Value = 199
$ws1.Range("F7").Value = 424
$ws1.Range("F8").Value = 421
$ws1.Range("F9").Value = 470
$ws1.Range("F10").Value = 260
$ws1.Range("F11").Value = 12759
$ws1.Range("F12").Value = 45
$ws1.Range("F13").Value = 321
$ws1.Range("F17").Value = 271
$ws1.Range("F19").Value = 191
$ws1.Range("F21").Value = 178
$ws1.Range("F22").Value = 2758
$ws1.Range("F29").Value = 1103
$ws1.Range("F30").Value = 4265
$ws1.Range("F32").Value = 3798
$ws1.Range("F33").Value = 844
$ws1.Range("F36").Value = 73
$ws1.Range("F37").Value = 1366
$ws1.Range("F39").Value = 789
$ws1.Range("F40").Value = 47
$ws1.Range("F41").Value = 137
$ws1.Range("F42").Value = 503
$ws1.Range("F43").Value = 683
$ws1.Range("F44").Value = 77
$ws1.Range("F45").Value = 158
$ws1.Range("F46").Value = 292
$ws1.Range("F47").Value = 122
$ws1.Range("F48").Value = 160
$ws1.Range("F49").Value = 174

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 0
$ws2.Range("F11").Value = 35
$ws2.Range("F13").Value = 65
$ws2.Range("F16").Value = 189

# Sheet: 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 65

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 10207
$ws4.Range("F6").Value = 737
$ws4.Range("F8").Value = 199
$ws4.Range("F9").Value = 424
$ws4.Range("F10").Value = 421
$ws4.Range("F11").Value = 470
$ws4.Range("F12").Value = 260
$ws4.Range("F13").Value = 12759
$ws4.Range("F14").Value = 321
$ws4.Range("F15").Value = 65
$ws4.Range("F16").Value = 271
$ws4.Range("F19").Value = 191
$ws4.Range("F21").Value = 178
$ws4.Range("F22").Value = 2758
$ws4.Range("F29").Value = 1103
$ws4.Range("F30").Value = 4265
$ws4.Range("F31").Value = 3798
$ws4.Range("F32").Value = 844
$ws4.Range("F35").Value = 73
$ws4.Range("F36").Value = 1366
$ws4.Range("F38").Value = 789
$ws4.Range("F39").Value = 47
$ws4.Range("F40").Value = 137
$ws4.Range("F41").Value = 503
$ws4.Range("F43").Value = 683
$ws4.Range("F44").Value = 77
$ws4.Range("F45").Value = 158
$ws4.Range("F46").Value = 292
$ws4.Range("F47").Value = 122
$ws4.Range("F48").Value = 160
$ws4.Range("F49").Value = 174

Write-Output "Applied all F-column updates."
